$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new blank rows above the two pre-existing placeholder rows
# (old row 10 and old row 11), pushing the placeholders down to rows 13/14.
$ws.Range("A10:A12").EntireRow.Insert()

# --- strEmailAccount: swap the outlook account used for the demo robot ---
$ws.Range("B5").Value = "Jason.Savory@defra.gov.uk"

# --- New rows 6-9: Outlook folder configuration for email processing ---
$ws.Range("A6").Value = "Email to process folder"
$ws.Range("B6").Value = "Inbox\COMPANY INFO\For Processing"
$ws.Range("C6").Value = "The folder to look for the mail for processing"

$ws.Range("A7").Value = "Email in progress folder"
$ws.Range("B7").Value = "Inbox\COMPANY INFO\In Progress"
$ws.Range("C7").Value = "The folder to move for in progress"

$ws.Range("A8").Value = "Email complete folder"
$ws.Range("B8").Value = "Inbox\COMPANY INFO\Completed"
$ws.Range("C8").Value = "The folder to move for completed mail"

$ws.Range("A9").Value = "Email exception folder"
$ws.Range("B9").Value = "Inbox\COMPANY INFO\Exceptions"
$ws.Range("C9").Value = "The folder to move exceptioned mail"
# Row 9 used to be the (taller) AttachmentDirectory row; restore its default height.
$ws.Rows.Item(9).AutoFit()

# --- Rows 10-12 (freshly inserted): existing parameters shifted down ---
$ws.Range("A10").Value = "StartPopUpTitle"
$ws.Range("B10").Value = "Companies House Demo Robot"
$ws.Range("C10").Value = "The title of the start pop up."

$ws.Range("A11").Value = "WorkpackageName"
$ws.Range("B11").Value = "Demo Robot"
$ws.Range("C11").Value = "workpackage name"

$ws.Range("A12").Value = "AttachmentDirectory"
$ws.Range("B12").Value = "{0}\Desktop\Attachments"
$ws.Range("C12").Value = "The file path for storing email attachments"
$ws.Rows.Item(12).RowHeight = 30

# --- Rows 13-14 (previously the blank placeholder rows): new parameters
#     for the OCRDeepDive / orchestrator integration work ---
$ws.Range("A13").Value = "RunTypeAssetName"
# Clear the left-over bold/italic/underline placeholder formatting so the
# cell renders like the other plain data cells.
$ws.Range("A13").Font.Bold = $false
$ws.Range("A13").Font.Italic = $false
$ws.Range("A13").Font.Underline = $false
$ws.Range("B13").Value = "DEMO_RUN_TYPE"
$ws.Range("C13").Value = "The asset name to control attended/unattended run"
$ws.Rows.Item(13).RowHeight = 30

$ws.Range("A14").Value = "Orchestrator Folder Path"
$ws.Range("B14").Value = "EA Root/Utilities"
$ws.Range("C14").Value = "The path to the orcehstrator folder containing this process"
$ws.Rows.Item(14).RowHeight = 30

# Grow Table1 to cover the new rows
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:C14"))

# Match the workbook's saved selection
$ws.Range("B14").Select()
